# "map builder.xlsx" - room list touch-up
#
# The sheet is a grid of dungeon room-type names. This change:
#   1. Renames the F3 cell from "empty_hallway" to "empty_passageway".
#      ("empty_passageway" already exists elsewhere in the sheet, and
#      "empty_hallway" was only ever used here, so once the text is
#      updated the old shared string becomes orphaned and is dropped
#      from the workbook's shared-string table on save - that's what
#      shifts all the higher <v> indices used by other cells down by
#      one; we don't need to touch those cells ourselves.)
#   2. Clears the stray formatted-but-empty G6 cell.
#   3. Moves the active selection from J9 to F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "empty_passageway"

$ws.Range("G6").Clear()

$ws.Range("F3").Select()
